{"js": "// Replace each two-digit multiplication equation in the table with its\n// new value. Every \"old\" string below is unique in the document, so a\n// plain (non-wildcard) search-and-replace is unambiguous.\nconst replacements = [\n  [\"62\u00d711=682\", \"22\u00d777=1694\"],\n  [\"32\u00d735=1120\", \"19\u00d797=1843\"],\n  [\"38\u00d780=3040\", \"42\u00d726=1092\"],\n  [\"43\u00d749=2107\", \"66\u00d755=3630\"],\n  [\"29\u00d739=1131\", \"99\u00d762=6138\"],\n  [\"18\u00d745=810\", \"18\u00d784=1512\"],\n  [\"14\u00d798=1372\", \"40\u00d797=3880\"],\n  [\"33\u00d729=957\", \"53\u00d760=3180\"],\n  [\"83\u00d740=3320\", \"41\u00d745=1845\"],\n  [\"29\u00d788=2552\", \"45\u00d717=765\"],\n  [\"95\u00d797=9215\", \"75\u00d734=2550\"],\n  [\"73\u00d753=3869\", \"28\u00d716=448\"],\n  [\"75\u00d735=2625\", \"61\u00d744=2684\"],\n  [\"61\u00d740=2440\", \"74\u00d765=4810\"],\n  [\"51\u00d768=3468\", \"35\u00d798=3430\"],\n  [\"96\u00d769=6624\", \"72\u00d771=5112\"],\n  [\"80\u00d775=6000\", \"12\u00d742=504\"],\n  [\"34\u00d767=2278\", \"65\u00d741=2665\"],\n  [\"46\u00d727=1242\", \"30\u00d764=1920\"],\n  [\"91\u00d791=8281\", \"49\u00d725=1225\"],\n  [\"82\u00d792=7544\", \"69\u00d788=6072\"],\n  [\"26\u00d726=676\", \"18\u00d734=612\"],\n  [\"72\u00d756=4032\", \"67\u00d782=5494\"],\n  [\"67\u00d737=2479\", \"98\u00d760=5880\"],\n  [\"47\u00d764=3008\", \"95\u00d765=6175\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication equation in the table with its\n# new value. Every \"old\" string is unique in the document, so Find/Replace\n# with MatchWholeWord off but MatchCase on is unambiguous per item.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"62\u00d711=682\", \"22\u00d777=1694\"),\n    @(\"32\u00d735=1120\", \"19\u00d797=1843\"),\n    @(\"38\u00d780=3040\", \"42\u00d726=1092\"),\n    @(\"43\u00d749=2107\", \"66\u00d755=3630\"),\n    @(\"29\u00d739=1131\", \"99\u00d762=6138\"),\n    @(\"18\u00d745=810\", \"18\u00d784=1512\"),\n    @(\"14\u00d798=1372\", \"40\u00d797=3880\"),\n    @(\"33\u00d729=957\", \"53\u00d760=3180\"),\n    @(\"83\u00d740=3320\", \"41\u00d745=1845\"),\n    @(\"29\u00d788=2552\", \"45\u00d717=765\"),\n    @(\"95\u00d797=9215\", \"75\u00d734=2550\"),\n    @(\"73\u00d753=3869\", \"28\u00d716=448\"),\n    @(\"75\u00d735=2625\", \"61\u00d744=2684\"),\n    @(\"61\u00d740=2440\", \"74\u00d765=4810\"),\n    @(\"51\u00d768=3468\", \"35\u00d798=3430\"),\n    @(\"96\u00d769=6624\", \"72\u00d771=5112\"),\n    @(\"80\u00d775=6000\", \"12\u00d742=504\"),\n    @(\"34\u00d767=2278\", \"65\u00d741=2665\"),\n    @(\"46\u00d727=1242\", \"30\u00d764=1920\"),\n    @(\"91\u00d791=8281\", \"49\u00d725=1225\"),\n    @(\"82\u00d792=7544\", \"69\u00d788=6072\"),\n    @(\"26\u00d726=676\", \"18\u00d734=612\"),\n    @(\"72\u00d756=4032\", \"67\u00d782=5494\"),\n    @(\"67\u00d737=2479\", \"98\u00d760=5880\"),\n    @(\"47\u00d764=3008\", \"95\u00d765=6175\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\nWrite-Output \"done\"\n"}
